$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2292.7693
$ws.Range("I62").Value = 2180
$ws.Range("J62").Value = 2668.6667
$ws.Range("K62").Value = 2180
$ws.Range("L62").Value = 2668.6667
$ws.Range("M62").Value = -1556
$ws.Range("N62").Value = -3916.6667
$ws.Range("H64").Value = 36756.266
$ws.Range("I64").Value = 3155.5
$ws.Range("J64").Value = 59156.777
$ws.Range("K64").Value = 3155.5
$ws.Range("L64").Value = 59156.777
$ws.Range("M64").Value = -2907.5
$ws.Range("N64").Value = -59652.777
$ws.Range("H65").Value = 2292.7693
$ws.Range("I65").Value = 2180
$ws.Range("J65").Value = 2668.6667
$ws.Range("K65").Value = 10900
$ws.Range("L65").Value = 13343.3335
$ws.Range("M65").Value = -7780
$ws.Range("N65").Value = -19583.3335
$ws.Range("H67").Value = 36756.266
$ws.Range("I67").Value = 3155.5
$ws.Range("J67").Value = 59156.777
$ws.Range("K67").Value = 3155.5
$ws.Range("L67").Value = 59156.777
$ws.Range("M67").Value = -2297.5
$ws.Range("N67").Value = -60872.777
$ws.Range("H112").Value = 1497.4286
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1497.4286
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 4492.2858
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -6708.2858
$ws.Range("H113").Value = 2738
$ws.Range("I113").Value = 2565.2632
$ws.Range("K113").Value = 2565.2632
$ws.Range("M113").Value = 688.7368000000001
$ws.Range("H137").Value = 2580.982
$ws.Range("I137").Value = 2875.7097
$ws.Range("J137").Value = 2200.2917
$ws.Range("K137").Value = 8627.1291
$ws.Range("L137").Value = 6600.875100000001
$ws.Range("M137").Value = -6077.1291
$ws.Range("N137").Value = -11700.8751

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 667.4666999999999
$ws.Range("I2").Value = 643.7143
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 643.7143
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -530.7143
$ws.Range("N2").Value = -1226
$ws.Range("H32").Value = 6083.906
$ws.Range("I32").Value = 5448.94
$ws.Range("K32").Value = 5448.94
$ws.Range("M32").Value = -5161.94
$ws.Range("H45").Value = 1842.5
$ws.Range("I45").Value = 1297
$ws.Range("J45").Value = 2824.4
$ws.Range("K45").Value = 1297
$ws.Range("L45").Value = 2824.4
$ws.Range("M45").Value = -920
$ws.Range("N45").Value = -3578.4
$ws.Range("H116").Value = 667.4666999999999
$ws.Range("I116").Value = 643.7143
$ws.Range("J116").Value = 1000
$ws.Range("K116").Value = 643.7143
$ws.Range("L116").Value = 1000
$ws.Range("M116").Value = 1650.2857
$ws.Range("N116").Value = -5588
$ws.Range("H123").Value = 21800
$ws.Range("J123").Value = 21800
$ws.Range("L123").Value = 21800
$ws.Range("N123").Value = -31600
$ws.Range("H132").Value = 223904.84
$ws.Range("I132").Value = 4969.8057
$ws.Range("J132").Value = 1012071
$ws.Range("K132").Value = 14909.4171
$ws.Range("L132").Value = 3036213
$ws.Range("M132").Value = -12379.4171
$ws.Range("N132").Value = -3041273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 667.4666999999999
$ws.Range("I3").Value = 643.7143
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 643.7143
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -529.7143
$ws.Range("N3").Value = -1228
$ws.Range("H99").Value = 1541.4286
$ws.Range("I99").Value = 1782.8572
$ws.Range("J99").Value = 1300
$ws.Range("K99").Value = 1782.8572
$ws.Range("L99").Value = 1300
$ws.Range("M99").Value = -284.8571999999999
$ws.Range("N99").Value = -4296

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1393.5
$ws.Range("I16").Value = 1352.75
$ws.Range("K16").Value = 1352.75
$ws.Range("M16").Value = -1065.75
$ws.Range("H86").Value = 2610.7368
$ws.Range("I86").Value = 2391.3333
$ws.Range("J86").Value = 2986.8572
$ws.Range("K86").Value = 2391.3333
$ws.Range("L86").Value = 2986.8572
$ws.Range("M86").Value = -1268.3333
$ws.Range("N86").Value = -5232.8572
$ws.Range("H89").Value = 2610.7368
$ws.Range("I89").Value = 2391.3333
$ws.Range("J89").Value = 2986.8572
$ws.Range("K89").Value = 11956.6665
$ws.Range("L89").Value = 14934.286
$ws.Range("M89").Value = -6340.666499999999
$ws.Range("N89").Value = -26166.286
$ws.Range("H99").Value = 1239.7084
$ws.Range("I99").Value = 1100.579
$ws.Range("J99").Value = 1768.4
$ws.Range("K99").Value = 1100.579
$ws.Range("L99").Value = 1768.4
$ws.Range("M99").Value = 397.421
$ws.Range("N99").Value = -4764.4
$ws.Range("H105").Value = 1393.5
$ws.Range("I105").Value = 1301
$ws.Range("J105").Value = 1560
$ws.Range("K105").Value = 1301
$ws.Range("L105").Value = 1560
$ws.Range("M105").Value = 446
$ws.Range("N105").Value = -5054
$ws.Range("H107").Value = 825.16
$ws.Range("I107").Value = 812.5789
$ws.Range("J107").Value = 865
$ws.Range("K107").Value = 812.5789
$ws.Range("L107").Value = 865
$ws.Range("M107").Value = 1107.4211
$ws.Range("N107").Value = -4705
$ws.Range("H113").Value = 1393.5
$ws.Range("I113").Value = 1352.75
$ws.Range("K113").Value = 1352.75
$ws.Range("M113").Value = 817.25
$ws.Range("H122").Value = 3275.7273
$ws.Range("I122").Value = 2742.6
$ws.Range("J122").Value = 3720
$ws.Range("K122").Value = 8227.799999999999
$ws.Range("L122").Value = 11160
$ws.Range("M122").Value = -5777.799999999999
$ws.Range("N122").Value = -16060
$ws.Range("H126").Value = 1239.7084
$ws.Range("I126").Value = 1100.579
$ws.Range("J126").Value = 1768.4
$ws.Range("K126").Value = 3301.737
$ws.Range("L126").Value = 5305.200000000001
$ws.Range("M126").Value = -831.7370000000001
$ws.Range("N126").Value = -10245.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 9487.241
$ws.Range("I132").Value = 6580.5
$ws.Range("J132").Value = 23439.6
$ws.Range("K132").Value = 19741.5
$ws.Range("L132").Value = 70318.79999999999
$ws.Range("M132").Value = -17211.5
$ws.Range("N132").Value = -75378.79999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1772.6666
$ws.Range("I7").Value = 1637.7693
$ws.Range("J7").Value = 1932.091
$ws.Range("K7").Value = 1637.7693
$ws.Range("L7").Value = 1932.091
$ws.Range("M7").Value = -1525.7693
$ws.Range("N7").Value = -2156.091
$ws.Range("H82").Value = 2211.5
$ws.Range("I82").Value = 1081
$ws.Range("J82").Value = 2646.3076
$ws.Range("K82").Value = 1081
$ws.Range("L82").Value = 2646.3076
$ws.Range("M82").Value = -720
$ws.Range("N82").Value = -3368.3076
$ws.Range("H85").Value = 2211.5
$ws.Range("I85").Value = 1081
$ws.Range("J85").Value = 2646.3076
$ws.Range("K85").Value = 1081
$ws.Range("L85").Value = 2646.3076
$ws.Range("M85").Value = 167
$ws.Range("N85").Value = -5142.3076
$ws.Range("H93").Value = 1787.75
$ws.Range("I93").Value = 1464
$ws.Range("J93").Value = 2500
$ws.Range("K93").Value = 1464
$ws.Range("L93").Value = 2500
$ws.Range("M93").Value = -216
$ws.Range("N93").Value = -4996
$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H126").Value = 1772.6666
$ws.Range("I126").Value = 1637.7693
$ws.Range("J126").Value = 1932.091
$ws.Range("K126").Value = 4913.3079
$ws.Range("L126").Value = 5796.272999999999
$ws.Range("M126").Value = -2443.3079
$ws.Range("N126").Value = -10736.273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 531.0909
$ws.Range("I107").Value = 362.85715
$ws.Range("J107").Value = 825.5
$ws.Range("K107").Value = 1088.57145
$ws.Range("L107").Value = 2476.5
$ws.Range("M107").Value = 831.4285500000001
$ws.Range("N107").Value = -6316.5
$ws.Range("H113").Value = 657.2727
$ws.Range("J113").Value = 758
$ws.Range("L113").Value = 2274
$ws.Range("N113").Value = -6614
$ws.Range("H122").Value = 668609.75
$ws.Range("I122").Value = 910876.9399999999
$ws.Range("J122").Value = 2375
$ws.Range("K122").Value = 2732630.82
$ws.Range("L122").Value = 7125
$ws.Range("M122").Value = -2730180.82
$ws.Range("N122").Value = -12025
